# Aggiornamento fino a 6/03: append nuovi dati per Carpi (righe 245-247)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Imposta i valori delle nuove righe (date come seriali, in linea con i dati esistenti)
$ws.Cells.Item(245, 1).Value = 44319
$ws.Cells.Item(245, 2).Value = 11
$ws.Cells.Item(245, 3).Value = 98
$ws.Cells.Item(245, 4).Value = 135.417098481394

$ws.Cells.Item(246, 1).Value = 44320
$ws.Cells.Item(246, 2).Value = 4
$ws.Cells.Item(246, 3).Value = 87
$ws.Cells.Item(246, 4).Value = 120.2172200804212

$ws.Cells.Item(247, 1).Value = 44321
$ws.Cells.Item(247, 2).Value = 1
$ws.Cells.Item(247, 3).Value = 86
$ws.Cells.Item(247, 4).Value = 118.83541295306

# Copia la formattazione della colonna data (colonna A) dall'ultima riga esistente
# alle nuove righe, cosi' da mantenere lo stesso stile (grassetto, bordo, centrato,
# formato data) gia' usato per tutte le altre righe della colonna A.
$ws.Cells.Item(244, 1).Copy()
$ws.Cells.Item(245, 1).PasteSpecial(-4122)
$ws.Cells.Item(246, 1).PasteSpecial(-4122)
$ws.Cells.Item(247, 1).PasteSpecial(-4122)
